# Update cryptos list (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold plain text (e.g. "1.000", "0.02000") that
# Excel would otherwise auto-convert to numbers (stripping significant
# trailing/leading zeros). Force them to stay text before writing values.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("D2").Value2 = "27.208.70"
$ws.Range("E2").Value2 = "  +0.52%  "

# --- Row 3 ---
$ws.Range("D3").Value2 = "1.903.32"
$ws.Range("E3").Value2 = "  +0.51%  "

# --- Row 4 ---
$ws.Range("D4").Value2 = "1.000"
$ws.Range("E4").Value2 = "  -0.26%  "

# --- Row 5 ---
$ws.Range("D5").Value2 = "306.01"
$ws.Range("E5").Value2 = "  -0.18%  "

# --- Row 6 ---
$ws.Range("D6").Value2 = "1.000"
$ws.Range("E6").Value2 = "  -0.21%  "

# --- Row 7 ---
$ws.Range("D7").Value2 = "0.5415"
$ws.Range("E7").Value2 = "  +3.87%  "

# --- Row 8 ---
$ws.Range("D8").Value2 = "0.3804"
$ws.Range("E8").Value2 = "  +1.20%  "

# --- Row 9 ---
$ws.Range("D9").Value2 = "0.07288"
$ws.Range("E9").Value2 = "  +0.44%  "

# --- Row 10 ---
$ws.Range("D10").Value2 = "22.08"
$ws.Range("E10").Value2 = "  +4.72%  "

# --- Row 11 ---
$ws.Range("D11").Value2 = "0.9026"
$ws.Range("E11").Value2 = "  +0.45%  "

# --- Row 12 ---
$ws.Range("D12").Value2 = "0.08185"
$ws.Range("E12").Value2 = "  +0.10%  "

# --- Row 13 (only price changed) ---
$ws.Range("D13").Value2 = "95.39"

# --- Row 14 ---
$ws.Range("D14").Value2 = "5.340"
$ws.Range("E14").Value2 = "  +0.86%  "

# --- Row 15 ---
$ws.Range("D15").Value2 = "1.000"
$ws.Range("E15").Value2 = "  -0.31%  "

# --- Row 16 ---
$ws.Range("D16").Value2 = "14.80"
$ws.Range("E16").Value2 = "  +1.64%  "

# --- Row 17 (only price changed) ---
$ws.Range("D17").Value2 = "0.000008622"

# --- Row 18 (only volume changed) ---
$ws.Range("E18").Value2 = "  -0.31%  "

# --- Row 19 ---
$ws.Range("D19").Value2 = "1.310.76"
$ws.Range("E19").Value2 = "  -30.97%  "

# --- Row 20 ---
$ws.Range("D20").Value2 = "27.239.19"
$ws.Range("E20").Value2 = "  +0.51%  "

# --- Row 21 ---
$ws.Range("D21").Value2 = "5.043"
$ws.Range("E21").Value2 = "  -0.68%  "

# --- Row 23 ---
$ws.Range("D23").Value2 = "6.507"
$ws.Range("E23").Value2 = "  +1.60%  "

# --- Row 24 ---
$ws.Range("D24").Value2 = "148.24"
$ws.Range("E24").Value2 = "  -0.32%  "

# --- Row 25 ---
$ws.Range("D25").Value2 = "2.303"
$ws.Range("E25").Value2 = "  +0.63%  "

# --- Row 26 (only volume changed) ---
$ws.Range("E26").Value2 = "  +0.95%  "

# --- Row 27 ---
$ws.Range("D27").Value2 = "1.758"
$ws.Range("E27").Value2 = "  +1.23%  "

# --- Row 28 ---
$ws.Range("D28").Value2 = "116.60"
$ws.Range("E28").Value2 = "  +1.33%  "

# --- Row 29 ---
$ws.Range("D29").Value2 = "4.849"
$ws.Range("E29").Value2 = "  +1.38%  "

# --- Row 30 ---
$ws.Range("D30").Value2 = "4.646"
$ws.Range("E30").Value2 = "  -4.11%  "

# --- Row 31 ---
$ws.Range("D31").Value2 = "0.09201"
$ws.Range("E31").Value2 = "  -0.18%  "

# --- Row 32 ---
$ws.Range("D32").Value2 = "0.8211"
$ws.Range("E32").Value2 = "  +4.41%  "

# --- Row 33 ---
$ws.Range("D33").Value2 = "0.05062"
$ws.Range("E33").Value2 = "  +0.68%  "

# --- Row 34 ---
$ws.Range("D34").Value2 = "1.221"
$ws.Range("E34").Value2 = "  +0.92%  "

# --- Row 35 ---
$ws.Range("D35").Value2 = "3.008"
$ws.Range("E35").Value2 = "  +1.39%  "

# --- Row 36 ---
$ws.Range("D36").Value2 = "3.315"
$ws.Range("E36").Value2 = "  -3.11%  "

# --- Row 37 ---
$ws.Range("D37").Value2 = "2.695"
$ws.Range("E37").Value2 = "  +3.24%  "

# --- Row 38 ---
$ws.Range("D38").Value2 = "0.6024"
$ws.Range("E38").Value2 = "  +5.66%  "

# --- Row 39 ---
$ws.Range("D39").Value2 = "0.02000"
$ws.Range("E39").Value2 = "  +0.59%  "

# --- Row 40 ---
$ws.Range("D40").Value2 = "1.076"
$ws.Range("E40").Value2 = "  +0.12%  "

# --- Row 41 ---
$ws.Range("D41").Value2 = "9.257"
$ws.Range("E41").Value2 = "  +2.81%  "

# --- Row 42 ---
$ws.Range("D42").Value2 = "6.650"
$ws.Range("E42").Value2 = "  +1.54%  "

# --- Rows 43 & 44 swap places: Quant <-> Decentraland ---
$ws.Range("B43").Value2 = "Decentraland"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D43").Value2 = "0.5168"
$ws.Range("E43").Value2 = "  +6.60%  "

$ws.Range("B44").Value2 = "Quant"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value2 = "115.70"
$ws.Range("E44").Value2 = "  -0.24%  "

# --- Row 45 ---
$ws.Range("D45").Value2 = "0.1530"
$ws.Range("E45").Value2 = "  +1.12%  "

# --- Rows 46 & 47 swap places: EnergySwap <-> PaxDollar ---
$ws.Range("B46").Value2 = "PaxDollar"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value2 = "0.9997"
$ws.Range("E46").Value2 = "  -0.29%  "

$ws.Range("B47").Value2 = "EnergySwap"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value2 = "10.14"
$ws.Range("E47").Value2 = "  +0.54%  "

# --- Row 48 (only volume changed) ---
$ws.Range("E48").Value2 = "  +1.07%  "

# --- Row 49 ---
$ws.Range("D49").Value2 = "38.11"
$ws.Range("E49").Value2 = "  -0.13%  "

# --- Row 50 (only volume changed) ---
$ws.Range("E50").Value2 = "  +2.84%  "

# --- Row 51 ---
$ws.Range("D51").Value2 = "63.45"
$ws.Range("E51").Value2 = "  +0.02%  "
